$p = $ppt.ActivePresentation

# ----------------------------------------------------------------------
# Slide 2 ("AIR BATTLE 2014" / "Game project" title slide)
# ----------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$title = $s2.Shapes.Item(1)
$tr = $title.TextFrame.TextRange

# Replace "AIR BATTLE 2014" with "2015 THE AIR BATTLE" and italicize it.
$oldHeadLen = "AIR BATTLE 2014".Length
$tr.Characters(1, $oldHeadLen).Text = "2015 THE AIR BATTLE"
$newHeadLen = "2015 THE AIR BATTLE".Length
$tr.Characters(1, $newHeadLen).Font.Italic = $true

# "Game project" (after the line break) becomes "JS game", typed as two
# separate edits ("Game" -> "JS g", then " project" -> "ame") so the
# run ends up split the same way a live edit would split it.
$secondLineStart = $newHeadLen + 2
$tr.Characters($secondLineStart, 4).Text = "JS g"
$tr.Characters($secondLineStart + 4, 8).Text = "ame"

# ----------------------------------------------------------------------
# Slide 3 ("What we use:" content slide)
# ----------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$body = $s3.Shapes.Item(2)
$tr2 = $body.TextFrame.TextRange

$cr = [char]13
$newText = "KineticJS for HTML5 Canvas" + $cr + "Raphael for SVG graphics" + $cr + "Canvas and SVG animation" + $cr + "Repository on GitHub" + $cr + "JS OOP"
$tr2.Text = $newText

# Re-create the finer run splits that the reordered/retyped bullets have
# in the authored deck (same trick as above: replacing a sub-range with
# itself forces a run boundary without altering the visible text).
$tr2.Characters(1, 20).Text = "KineticJS for HTML5 "
$tr2.Characters(28, 16).Text = "Raphael for SVG "
$tr2.Characters(99, 3).Text = "JS "
